$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("GSMs")
$ws2 = $wb.Worksheets.Item("Values")

# --- Sheet "GSMs": highlight the mandatory-field headers (CallStatus, Region, Segment)
# by copying the red "GSM" header format onto C1, E1, G1.
$ws1.Range("A1").Copy()
$ws1.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("E1").PasteSpecial(-4122)
$ws1.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Sheet "Values": add a new "Segment" code column (C) with a single-letter
# code per CallStatus row, mirroring column A's header style / column B's data style.
$ws2.Range("C1").Value = "Segment"

# Match formatting: header like A1/B1 (style 2), data cells like B2:B21 (style 1).
$ws2.Range("B1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("B2:B22").Copy()
$ws2.Range("C2:C22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the one-letter codes (rows entered in the same order they were
# originally authored, row -> letter).
$order = @(
    @(2, "A"), @(6, "E"), @(17, "P"), @(20, "S"), @(19, "R"), @(7, "F"), @(13, "L"), @(14, "M"),
    @(4, "C"), @(3, "B"), @(5, "D"), @(8, "G"), @(9, "H"), @(10, "I"), @(11, "J"), @(12, "K"),
    @(15, "N"), @(16, "O"), @(18, "Q"), @(21, "T")
)
foreach ($pair in $order) {
    $ws2.Cells.Item($pair[0], 3).Value = $pair[1]
}
$ws2.Cells.Item(22, 3).ClearContents()

# Column width for the new column C (bestFit-style width matching sheet "GSMs" column C).
$ws2.Columns.Item(3).ColumnWidth = 8.83

# Zoom the Values sheet view.
$ws2.Activate()
$excel.ActiveWindow.Zoom = 85

# Selections (cosmetic, mirrors where the author last clicked).
[void]$ws2.Range("A1:C1").Select()
$ws1.Activate()
[void]$ws1.Range("I11").Select()
